$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (just below the fixed
# row 483), pushing the existing rows 484:536 down to 486:538 and growing
# the used range from A1:T536 to A1:T538.
$ws.Rows("484:485").Insert()

# New row 484: latest "1a amarillo" quote for Vega Monumental Concepción.
$ws.Cells.Item(484, 1).Value = 11
$ws.Cells.Item(484, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(484, 3).Value = "Bíobío"
$ws.Cells.Item(484, 4).Value = 44769
$ws.Cells.Item(484, 5).Value = 8
$ws.Cells.Item(484, 6).Value = "Fruta"
$ws.Cells.Item(484, 7).Value = 100102
$ws.Cells.Item(484, 8).Value = "Cítricos"
$ws.Cells.Item(484, 9).Value = 100102003
$ws.Cells.Item(484, 10).Value = "Limón"
$ws.Cells.Item(484, 11).Value = "Sin especificar"
$ws.Cells.Item(484, 12).Value = "1a amarillo"
$ws.Cells.Item(484, 13).Value = 350
$ws.Cells.Item(484, 14).Value = 4500
$ws.Cells.Item(484, 15).Value = 5000
$ws.Cells.Item(484, 16).Value = 4786
$ws.Cells.Item(484, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(484, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(484, 19).Value = 299
$ws.Cells.Item(484, 20).Value = 16

# New row 485: latest "2a amarillo" quote for Vega Monumental Concepción.
$ws.Cells.Item(485, 1).Value = 11
$ws.Cells.Item(485, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(485, 3).Value = "Bíobío"
$ws.Cells.Item(485, 4).Value = 44769
$ws.Cells.Item(485, 5).Value = 8
$ws.Cells.Item(485, 6).Value = "Fruta"
$ws.Cells.Item(485, 7).Value = 100102
$ws.Cells.Item(485, 8).Value = "Cítricos"
$ws.Cells.Item(485, 9).Value = 100102003
$ws.Cells.Item(485, 10).Value = "Limón"
$ws.Cells.Item(485, 11).Value = "Sin especificar"
$ws.Cells.Item(485, 12).Value = "2a amarillo"
$ws.Cells.Item(485, 13).Value = 120
$ws.Cells.Item(485, 14).Value = 4000
$ws.Cells.Item(485, 15).Value = 4000
$ws.Cells.Item(485, 16).Value = 4000
$ws.Cells.Item(485, 17).Value = "$/malla 16 kilos"
$ws.Cells.Item(485, 18).Value = "Región de O'Higgins"
$ws.Cells.Item(485, 19).Value = 250
$ws.Cells.Item(485, 20).Value = 16
